$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "item_num" column just before the existing "comment"
#    column (which is column O) on every sheet that has it (Measures, ID,
#    Dems, Dates, NewVars). The old column O data is shifted right to
#    column P, and the new column O gets the "item_num" header.
# ---------------------------------------------------------------------------

$sheetNames = @("Measures", "ID", "Dems", "Dates", "NewVars")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $lastRow = $ws.UsedRange.Rows.Count

    # Move the existing column O ("comment") to column P.
    $ws.Range("P1:P" + $lastRow).Value2 = $ws.Range("O1:O" + $lastRow).Value2
    $ws.Range("O1:O" + $lastRow).ClearContents()

    # New header for column O.
    $ws.Range("O1").Value2 = "item_num"

    # Update the view: selection moves from the old "comment" column
    # reference to the new one, matching the other unaffected columns.
    $ws.Activate()
    $ws.Range("O1:O1048576").Select() | Out-Null
}

# Measures sheet: fill in the item number for every data row (2-59) with 1.
$measures = $wb.Worksheets.Item("Measures")
$measures.Range("O2:O59").Value2 = 1

# ---------------------------------------------------------------------------
# 2. Filter database defined name must grow from column O to column P to
#    keep covering the full header row.
# ---------------------------------------------------------------------------

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Measures!_FilterDatabase") {
        $n.RefersTo = "=Measures!`$A`$1:`$P`$1"
    }
}

# ---------------------------------------------------------------------------
# 3. The active sheet moves from PanelInfo to NewVars.
# ---------------------------------------------------------------------------

$wb.Worksheets.Item("NewVars").Activate()
